$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Referensi" (sheet1): append 3 new reference rows (5-7)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Referensi")

$ws1.Range("A5").Value = "Free SFX"
$ws1.Range("B5").Value = "https://freesfx.co.uk/"
$ws1.Range("C5").Value = "Free SFX library"
$ws1.Range("D5").Value = "Pustaka efek suara gratis"

$ws1.Range("A6").Value = "cegaton"
$ws1.Range("B6").Value = "https://blender.stackexchange.com/questions/10725/how-do-i-create-an-equilateral-tetrahedron"
$ws1.Range("C6").Value = "Tetrahedron Blender"
$ws1.Range("D6").Value = "Buat dengan fungsi matematika"

$ws1.Range("A7").Value = "Robin bets"
$ws1.Range("B7").Value = "https://blender.stackexchange.com/questions/10725/how-do-i-create-an-equilateral-tetrahedron"
$ws1.Range("C7").Value = "Tetrahedron Blender"
$ws1.Range("D7").Value = "buat dengan sudut berlawanan"

# ---------------------------------------------------------------------------
# Sheet "Aset" (sheet2): append 3 new asset rows (5-7) + widen columns
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Aset")

$ws2.Range("A5").Value = "mynamesisjacobj"
$ws2.Range("B5").Value = "https://gist.github.com/mynameisjacobj/8584bdc1e9e104b005044fdfda5fe9e5"
$ws2.Range("C5").Value = "Non Convex Mesh Collider"
$ws2.Range("D5").Value = "Unity tidak mendukun mesh collider non-cembung"
$ws2.Range("E5").Value = "[<100%]"

$ws2.Range("A6").Value = "Productivity Boost"
$ws2.Range("B6").Value = "https://assetstore.unity.com/packages/tools/physics/non-convex-mesh-collider-84867"
$ws2.Range("C6").Value = "[ORI] Non Convex Mesh Collider"
$ws2.Range("D6").Value = "[Tidak digunakan, tidak gratis] buat Mesh Collider cekung!"
$ws2.Range("E6").Value = "[>100%]"

$ws2.Range("A7").Value = "Unity Asset Collection "
$ws2.Range("B7").Value = "http://unityassetcollection.com/non-convex-mesh-collider-free-download/"
$ws2.Range("C7").Value = "[Bajakan] Non Convex Mesh Collider"
$ws2.Range("D7").Value = "[Tidak digunakan, larangan, dilarang bajakan] buat Mesh Collider cekung!"
$ws2.Range("E7").Value = "[>100%]"

$ws2.Columns.Item(1).ColumnWidth = 29.88671875
$ws2.Columns.Item(2).ColumnWidth = 73.21875
$ws2.Columns.Item(3).ColumnWidth = 58.44140625
$ws2.Columns.Item(4).ColumnWidth = 62.33203125

# ---------------------------------------------------------------------------
# New sheet "Alat" (tools), inserted after "Aset"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$ws3.Name = "Alat"

$ws3.Range("A1").Value = "Penulis"
$ws3.Range("B1").Value = "Sumber Tautan"
$ws3.Range("C1").Value = "Judul"
$ws3.Range("D1").Value = "Deskripsi"
$ws3.Range("E1").Value = "Hak Cipta"

$ws3.Range("A2").Value = "Unity Technologies"
$ws3.Range("B2").Value = "https://www.unity3d.com"
$ws3.Hyperlinks.Add($ws3.Range("B2"), "https://www.unity3d.com", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "https://www.unity3d.com")
$ws3.Range("B2").Style = "Hipertaut"
$ws3.Range("C2").Value = "Unity Engine"
$ws3.Range("D2").Value = "buat game 2D atau 3D"
$ws3.Range("E2").Value = "[100%]"

$ws3.Range("A3").Value = "Blender Foundation"
$ws3.Range("B3").Value = "https://www.blender.org"
$ws3.Hyperlinks.Add($ws3.Range("B3"), "https://www.blender.org", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "https://www.blender.org")
$ws3.Range("B3").Style = "Hipertaut"
$ws3.Range("C3").Value = "Blender 3D"
$ws3.Range("D3").Value = "Buat apapun 3D"
$ws3.Range("E3").Value = "[<100%]"

$ws3.Columns.Item(1).ColumnWidth = 29.88671875
$ws3.Columns.Item(2).ColumnWidth = 73.21875
$ws3.Columns.Item(3).ColumnWidth = 58.44140625
$ws3.Columns.Item(4).ColumnWidth = 62.33203125
$ws3.Columns.Item(5).ColumnWidth = 20.21875

# ---------------------------------------------------------------------------
# Selection / active-view bookkeeping to mirror the target view state
# ---------------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("E10").Select()

$ws3.Activate()
$ws3.Range("E2").Select()

$ws1.Activate()
$ws1.Range("D7").Select()
